# Add two new columns, I (header "I0") and J (header "IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - same bold/border/centered style as the other
# header cells (e.g. column H).
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for rows 2-15, columns I (9) and J (10).
$data = @{
    2  = @(7, 8)
    3  = @(10, 10)
    4  = @(4, 6)
    5  = @(7, 9)
    6  = @(6, 7)
    7  = @(1, 1)
    8  = @(7, 7)
    9  = @(8, 8)
    10 = @(4, 4)
    11 = @(5, 5)
    12 = @(4, 4)
    13 = @(5, 6)
    14 = @(7, 7)
    15 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
